$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "29.533.69"
$ws.Range("E2").Value = "  +1.95%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.842.62"
$ws.Range("E3").Value = "  +0.92%  "

# Row 4 - TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9984"
$ws.Range("E4").Value = "  -0.61%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.76"
$ws.Range("E5").Value = "  +1.06%  "

# Row 6 - XRP
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6288"
$ws.Range("E6").Value = "  +1.60%  "

# Row 7 - USDC
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9992"
$ws.Range("E7").Value = "  -0.62%  "

# Row 8 - Dogecoin
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07460"
$ws.Range("E8").Value = "  +1.07%  "

# Row 9 - Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2964"
$ws.Range("E9").Value = "  +2.12%  "

# Row 10 - Solana
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.75"
$ws.Range("E10").Value = "  +3.93%  "

# Row 11 - TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07668"
$ws.Range("E11").Value = "  -0.10%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "1.835.68"
$ws.Range("E12").Value = "  +0.34%  "

# Row 13 - Polkadot
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.034"
$ws.Range("E13").Value = "  +1.64%  "

# Row 14 - Polygon
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6801"
$ws.Range("E14").Value = "  +2.11%  "

# Row 15 - Litecoin
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "84.24"
$ws.Range("E15").Value = "  +2.73%  "

# Row 16 - ShibaInu
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009375"
$ws.Range("E16").Value = "  +3.73%  "

# Row 17 - Uniswap
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.948"
$ws.Range("E17").Value = "  +1.65%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "29.505.57"
$ws.Range("E18").Value = "  +1.78%  "

# Row 19 - WrappedliquidstakedEther2.0
$ws.Range("D19").Value = "2.078.95"
$ws.Range("E19").Value = "  -0.07%  "

# Row 20 - BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "237.69"
$ws.Range("E20").Value = "  +1.29%  "

# Row 21 - Avalanche
$ws.Range("E21").Value = "  +0.90%  "

# Row 22 - Dai
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9992"
$ws.Range("E22").Value = "  -0.66%  "

# Row 23 - Chainlink
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.381"
$ws.Range("E23").Value = "  +3.71%  "

# Row 24 - BinanceUSD
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9997"
$ws.Range("E24").Value = "  -0.66%  "

# Row 25 - Monero
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.36"
$ws.Range("E25").Value = "  +0.26%  "

# Row 26 - Stellar
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1417"
$ws.Range("E26").Value = "  +0.72%  "

# Row 27 - Cosmos
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.542"
$ws.Range("E27").Value = "  +0.95%  "

# Row 28 - EthereumClassic
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.82"
$ws.Range("E28").Value = "  +0.73%  "

# Row 29 - Hedera
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06100"
$ws.Range("E29").Value = "  +10.30%  "

# Row 30 - PancakeSwap
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.499"
$ws.Range("E30").Value = "  +0.39%  "

# Row 31 - Toncoin
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.243"
$ws.Range("E31").Value = "  +2.76%  "

# Row 32 - Filecoin(was ICP)
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.127"
$ws.Range("E32").Value = "  +0.75%  "

# Row 33 - ICP(was Filecoin)
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.106"
$ws.Range("E33").Value = "  +0.15%  "

# Row 34 - LidoDAOToken
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.871"
$ws.Range("E34").Value = "  +2.19%  "

# Row 35 - ARBITRUM
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.147"
$ws.Range("E35").Value = "  +1.42%  "

# Row 36 - ImmutableX
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7297"
$ws.Range("E36").Value = "  -0.53%  "

# Row 37 - HuobiToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.605"
$ws.Range("E37").Value = "  -1.48%  "

# Row 38 - MXToken
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.877"
$ws.Range("E38").Value = "  +2.29%  "

# Row 39 - Maker
$ws.Range("D39").Value = "1.224.80"
$ws.Range("E39").Value = "  +2.29%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  +0.14%  "

# Row 41 - FraxShare
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.294"
$ws.Range("E41").Value = "  -1.12%  "

# Row 42 - TrustWalletToken
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9149"
$ws.Range("E42").Value = "  +1.18%  "

# Row 43 - PaxDollar
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  -0.47%  "

# Row 44 - Quant
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.05"
$ws.Range("E44").Value = "  +1.14%  "

# Row 45 - RocketPoolETH
$ws.Range("D45").Value = "1.993.69"
$ws.Range("E45").Value = "  +0.97%  "

# Row 46 - Aave
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.74"
$ws.Range("E46").Value = "  +2.25%  "

# Row 47 - Mantle(was BabyDogeCoin)
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5077"
$ws.Range("E47").Value = "  -0.53%  "

# Row 48 - EnergySwap(was Mantle)
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.235"
$ws.Range("E48").Value = "  +1.75%  "

# Row 49 - TheSandbox(was EnergySwap)
$ws.Range("B49").Value = "TheSandbox"
$ws.Range("C49").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4071"
$ws.Range("E49").Value = "  +1.66%  "

# Row 50 - BabyDogeCoin(was TheSandbox)
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00000000118"
$ws.Range("E50").Value = "  -1.86%  "

# Row 51 - Algorand
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1142"
$ws.Range("E51").Value = "  +4.11%  "
